$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2025-08-06 05:06:29"
$ws.Range("B14").Value = "create-repo"
$ws.Range("C14").Value = "new-organization97"
$ws.Range("H14").Value = "desk"

$i14 = $ws.Range("I14")
$i14.Value = "'False"
$i14.ClearFormats()
